$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- 1. Insert a new pinned row at the top of the table data (row 2) for the
#        new "Licht fixen!" task, shifting all existing rows down by one, then
#        grow the table to include the extra row. ---
$ws.Rows.Item(2).Insert()
$tbl.Resize($ws.Range("A1:E51"))

$ws.Range("A2").Value = "Licht fixen!"
$ws.Range("C2").Value = 0

# --- 2. The row that used to hold "last food too late horror event" (now at
#        row 4 after the shift above) is replaced by a new, mostly-empty task
#        row describing the follow-up work. ---
$ws.Range("A4").Value = "improve last food too late horror event"
$ws.Range("B4").Value = "Unclear"
$ws.Range("C4:E4").ClearContents()

# --- 3. Italicize the "Glas und Fensterscherben..." task (now row 3). ---
$ws.Range("A3").Font.Italic = $true

# --- 4. Extend the phone-horror-event note with the extra follow-up line. ---
$ws.Range("E14").Value = '"Hello?", "Who''s there?", "Weird", "Wrong number I hope..."'

# --- 5. One more trailing blank row appears at the bottom of the sheet
#        (mirrors the existing blank rows 52-58 styled with s=2). ---
$ws.Range("E59").Style = $ws.Range("E58").Style

# --- 6. Selection moves to A2 after the edits. ---
$ws.Range("A2").Select()
